$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Traceability-Matrix: link Order History & Shopping Cart test cases (column C)
# to their WireFrame design IDs.

# Order History wireframes
$ws.Range("C34").Value = "WireFrame_OrderHistory_001"
$ws.Range("C35").Value = "WireFrame_OrderHistory_002"
$ws.Range("C36").Value = "WireFrame_OrderHistory_001"
$ws.Range("C37").Value = "WireFrame_OrderHistory_001"

# Shopping cart wireframes
$ws.Range("C40").Value = "WireFrame_Cart_001"
$ws.Range("C41").Value = "WireFrame_Cart_001"
$ws.Range("C42").Value = "WireFrame_shoppingcart_003"
$ws.Range("C43").Value = "WireFrame_shoppingcart_002"
$ws.Range("C44").Value = "WireFrame_shoppingcart_004"

# C46 already held a WireFrame value (Wireframe_Checkout_002) but had no cell
# border/style applied; give it the same bordered look as the rest of the
# column by copying the formatting from a neighboring, already-styled cell.
$ws.Range("C44").Copy() | Out-Null
$ws.Range("C46").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Restore the last-saved selection
$ws.Range("C45").Select() | Out-Null
